$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J7: fill in date 2019-03-21 (serial 43545), matching the date style already
# used by the neighboring I7/K7 cells on that row.
$ws.Range("J7").Value = 43545
$ws.Range("J7").NumberFormat = $ws.Range("I7").NumberFormat

# I10: fill in date 2019-03-10 (serial 43534), matching K10's date style.
$ws.Range("I10").Value = 43534
$ws.Range("I10").NumberFormat = $ws.Range("K10").NumberFormat

# J10: fill in date 2019-03-13 (serial 43537), matching K10's date style.
$ws.Range("J10").Value = 43537
$ws.Range("J10").NumberFormat = $ws.Range("K10").NumberFormat

# Move the active selection from I5 to I14.
$ws.Range("I14").Select()
